$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: average_doctor / average_doctor_old columns were swapped (BP <-> BQ)

# Row 1
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# Row 4
$ws.Range("E4").Value = 0.428
$ws.Range("F4").Value = 0.07099999999999999
$ws.Range("G4").Value = 0.266
$ws.Range("N4").Value = 0.436
$ws.Range("O4").Value = 0.063
$ws.Range("P4").Value = 0.251
$ws.Range("Q4").Value = 0.024
$ws.Range("R4").Value = 0.017
$ws.Range("S4").Value = 0.129
$ws.Range("W4").Value = 0.286
$ws.Range("X4").Value = 0.109
$ws.Range("Y4").Value = 0.33
$ws.Range("AI4").Value = 0.292
$ws.Range("AJ4").Value = 0.08799999999999999
$ws.Range("AK4").Value = 0.297
$ws.Range("AU4").Value = 0.19
$ws.Range("AV4").Value = 0.028
$ws.Range("AW4").Value = 0.167
$ws.Range("BA4").Value = 2.001
$ws.Range("BB4").Value = 0.158
$ws.Range("BC4").Value = 0.397
$ws.Range("BG4").Value = 0.729
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.376
$ws.Range("BM4").Value = 0.716
$ws.Range("BN4").Value = 0.08
$ws.Range("BO4").Value = 0.282
$ws.Range("BP4").Value = 0.667
$ws.Range("BQ4").Value = 0.708

# Row 5
$ws.Range("E5").Value = 0.544
$ws.Range("F5").Value = 0.083
$ws.Range("G5").Value = 0.289
$ws.Range("N5").Value = 0.741
$ws.Range("O5").Value = 0.076
$ws.Range("P5").Value = 0.276
$ws.Range("Q5").Value = 0.016
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.083
$ws.Range("W5").Value = 0.276
$ws.Range("X5").Value = 0.109
$ws.Range("Y5").Value = 0.33
$ws.Range("AI5").Value = 0.312
$ws.Range("AJ5").Value = 0.099
$ws.Range("AK5").Value = 0.314
$ws.Range("AU5").Value = 0.369
$ws.Range("AV5").Value = 0.093
$ws.Range("AW5").Value = 0.305
$ws.Range("BA5").Value = 1.335
$ws.Range("BB5").Value = 0.079
$ws.Range("BC5").Value = 0.282
$ws.Range("BG5").Value = 0.397
$ws.Range("BH5").Value = 0.051
$ws.Range("BI5").Value = 0.226
$ws.Range("BM5").Value = 0.549
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.253
$ws.Range("BP5").Value = 0.445
$ws.Range("BQ5").Value = 0.457

# Row 6
$ws.Range("E6").Value = 0.479
$ws.Range("N6").Value = 0.549
$ws.Range("Q6").Value = 0.019
$ws.Range("W6").Value = 0.281
$ws.Range("AI6").Value = 0.302
$ws.Range("AU6").Value = 0.251
$ws.Range("BA6").Value = 1.593
$ws.Range("BG6").Value = 0.514
$ws.Range("BM6").Value = 0.621
$ws.Range("BP6").Value = 0.531
$ws.Range("BQ6").Value = 0.553

# Row 7
$ws.Range("E7").Value = 0.516
$ws.Range("N7").Value = 0.65
$ws.Range("Q7").Value = 0.017
$ws.Range("W7").Value = 0.278
$ws.Range("AI7").Value = 0.308
$ws.Range("AU7").Value = 0.31
$ws.Range("BA7").Value = 1.427
$ws.Range("BG7").Value = 0.437
$ws.Range("BM7").Value = 0.576
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.491

# Row 8
$ws.Range("E8").Value = 0.611
$ws.Range("F8").Value = 0.109
$ws.Range("G8").Value = 0.33
$ws.Range("N8").Value = 0.781
$ws.Range("O8").Value = 0.06
$ws.Range("P8").Value = 0.245
$ws.Range("Q8").Value = 0.017
$ws.Range("W8").Value = 0.304
$ws.Range("AI8").Value = 0.334
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.359
$ws.Range("AU8").Value = 0.311
$ws.Range("AV8").Value = 0.08400000000000001
$ws.Range("AW8").Value = 0.291
$ws.Range("BA8").Value = 1.741
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.352
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.107
$ws.Range("BI8").Value = 0.326
$ws.Range("BM8").Value = 0.6889999999999999
$ws.Range("BN8").Value = 0.067
$ws.Range("BO8").Value = 0.259
$ws.Range("BP8").Value = 0.58
$ws.Range("BQ8").Value = 0.604

# Row 9
$ws.Range("E9").Value = 0.5590000000000001
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.496
$ws.Range("N9").Value = 0.6879999999999999
$ws.Range("O9").Value = 0.215
$ws.Range("P9").Value = 0.463
$ws.Range("W9").Value = 0.204
$ws.Range("X9").Value = 0.163
$ws.Range("Y9").Value = 0.403
$ws.Range("AI9").Value = 0.258
$ws.Range("AJ9").Value = 0.191
$ws.Range("AK9").Value = 0.438
$ws.Range("BA9").Value = 1.699
$ws.Range("BB9").Value = 0.248
$ws.Range("BC9").Value = 0.498
$ws.Range("BG9").Value = 0.602
$ws.Range("BH9").Value = 0.24
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.645
$ws.Range("BN9").Value = 0.229
$ws.Range("BO9").Value = 0.478
$ws.Range("BP9").Value = 0.5659999999999999
$ws.Range("BQ9").Value = 0.584

# Row 10
$ws.Range("E10").Value = 0.6879999999999999
$ws.Range("F10").Value = 0.215
$ws.Range("G10").Value = 0.463
$ws.Range("N10").Value = 0.882
$ws.Range("O10").Value = 0.104
$ws.Range("P10").Value = 0.323
$ws.Range("W10").Value = 0.376
$ws.Range("X10").Value = 0.235
$ws.Range("Y10").Value = 0.484
$ws.Range("AI10").Value = 0.366
$ws.Range("AJ10").Value = 0.232
$ws.Range("AK10").Value = 0.482
$ws.Range("AU10").Value = 0.301
$ws.Range("AV10").Value = 0.21
$ws.Range("AW10").Value = 0.459
$ws.Range("BA10").Value = 2.076
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.656
$ws.Range("BH10").Value = 0.226
$ws.Range("BI10").Value = 0.475
$ws.Range("BM10").Value = 0.839
$ws.Range("BN10").Value = 0.135
$ws.Range("BO10").Value = 0.368
$ws.Range("BP10").Value = 0.6919999999999999
$ws.Range("BQ10").Value = 0.724

# Row 11
$ws.Range("E11").Value = 0.72
$ws.Range("F11").Value = 0.201
$ws.Range("G11").Value = 0.449
$ws.Range("N11").Value = 0.903
$ws.Range("O11").Value = 0.08699999999999999
$ws.Range("P11").Value = 0.296
$ws.Range("W11").Value = 0.376
$ws.Range("X11").Value = 0.235
$ws.Range("Y11").Value = 0.484
$ws.Range("AI11").Value = 0.398
$ws.Range("AJ11").Value = 0.24
$ws.Range("AK11").Value = 0.489
$ws.Range("AU11").Value = 0.441
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.496
$ws.Range("BA11").Value = 2.076
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.656
$ws.Range("BH11").Value = 0.226
$ws.Range("BI11").Value = 0.475
$ws.Range("BM11").Value = 0.839
$ws.Range("BN11").Value = 0.135
$ws.Range("BO11").Value = 0.368
$ws.Range("BP11").Value = 0.6919999999999999
$ws.Range("BQ11").Value = 0.727

# Row 12
$ws.Range("E12").Value = 1.403
$ws.Range("F12").Value = 0.748
$ws.Range("G12").Value = 0.865
$ws.Range("N12").Value = 1.465
$ws.Range("O12").Value = 1.039
$ws.Range("P12").Value = 1.02
$ws.Range("W12").Value = 1.629
$ws.Range("X12").Value = 0.576
$ws.Range("Y12").Value = 0.759
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.767
$ws.Range("AV12").Value = 2.737
$ws.Range("AW12").Value = 1.654
$ws.Range("BA12").Value = 3.708
$ws.Range("BB12").Value = 0.401
$ws.Range("BC12").Value = 0.633
$ws.Range("BG12").Value = 1.098
$ws.Range("BH12").Value = 0.121
$ws.Range("BI12").Value = 0.349
$ws.Range("BM12").Value = 1.295
$ws.Range("BN12").Value = 0.336
$ws.Range("BO12").Value = 0.58
$ws.Range("BP12").Value = 1.236
$ws.Range("BQ12").Value = 1.263

# Row 13
$ws.Range("E13").Value = 1.579
$ws.Range("F13").Value = 0.656
$ws.Range("G13").Value = 0.8100000000000001
$ws.Range("N13").Value = 2.063
$ws.Range("O13").Value = 0.9360000000000001
$ws.Range("P13").Value = 0.968
$ws.Range("W13").Value = 1.037
$ws.Range("X13").Value = 0.193
$ws.Range("Y13").Value = 0.439
$ws.Range("AI13").Value = 1.28
$ws.Range("AJ13").Value = 0.37
$ws.Range("AK13").Value = 0.608
$ws.Range("AU13").Value = 2.285
$ws.Range("AV13").Value = 0.925
$ws.Range("AW13").Value = 0.962
$ws.Range("BA13").Value = 2.353
$ws.Range("BB13").Value = 0.297
$ws.Range("BC13").Value = 0.545
$ws.Range("BG13").Value = 0.585
$ws.Range("BH13").Value = 0.07099999999999999
$ws.Range("BI13").Value = 0.267
$ws.Range("BM13").Value = 0.896
$ws.Range("BN13").Value = 0.283
$ws.Range("BO13").Value = 0.532
$ws.Range("BP13").Value = 0.784
$ws.Range("BQ13").Value = 0.73
